$d = $word.ActiveDocument

$r0 = $d.Content
$r0.Find.Execute("Video Title") | Out-Null
$r0.Text = "Titre de la vidéo"

$r1 = $d.Content
$r1.Find.Execute("Topic") | Out-Null
$r1.Text = "Rubrique"

$r2 = $d.Content
$r2.Find.Execute("Logic") | Out-Null
$r2.Text = "Logique"

$r3 = $d.Content
$r3.Find.Execute("Aim(s)") | Out-Null
$r3.Text = "Objectif(s)"

$r4 = $d.Content
$r4.Find.Execute("Train the ability to extrapolate information from a problem, work by exclusion.") | Out-Null
$r4.Text = "Former la capacité d’extrapoler des informations à partir d’un problème, de travailler par exclusion."

$r5 = $d.Content
$r5.Find.Execute("Length") | Out-Null
$r5.Text = "Durée"

$r6 = $d.Content
$r6.Find.Execute("Camp Location") | Out-Null
$r6.Text = "Lieu du camp"

$r7 = $d.Content
$r7.Find.Execute("Facilitators") | Out-Null
$r7.Text = "Animateurs"

$r8 = $d.Content
$r8.Find.Execute("N. of students") | Out-Null
$r8.Text = "N. des étudiants"

$r9 = $d.Content
$r9.Find.Execute("Resources") | Out-Null
$r9.Text = "Les ressources"

$r10 = $d.Content
$r10.Find.Execute("needed") | Out-Null
$r10.Text = "nécessaires"

$r11 = $d.Content
$r11.Find.Execute("Preparations") | Out-Null
$r11.Text = "Préparations"

$r12 = $d.Content
$r12.Find.Execute("Video time") | Out-Null
$r12.Text = "Temps de la vidéo"

$r13 = $d.Content
$r13.Find.Execute("What facilitator does") | Out-Null
$r13.Text = "Ce que fait le facilitateur"

$r14 = $d.Content
$r14.Find.Execute("What learners do") | Out-Null
$r14.Text = "Ce que font les apprenants"

$r15 = $d.Content
$r15.Find.Execute("General VMC Video Introduction") | Out-Null
$r15.Text = "Vidéo générale introduisant le CVM"

$r16 = $d.Content
$r16.Find.Execute("Video Introduction") | Out-Null
$r16.Text = "Video d'introduction"

$r17 = $d.Content
$r17.Find.Execute("00:47 – 02:11") | Out-Null
$r17.Text = "00:47 - 02:11"

$r18 = $d.Content
$r18.Find.Execute("Riddle") | Out-Null
$r18.Text = "Énigme"

$r19 = $d.Content
$r19.Find.Execute("Assist the process, provoke thoughts") | Out-Null
$r19.Text = "Faciliter le processus, susciter des pensées"

$r20 = $d.Content
$r20.Find.Execute("When a possible solution is suggested, ask the learners to enact the series of questions in the case of the suggested solution and check that each answer can be explained.") | Out-Null
$r20.Text = "Lorsqu'une solution possible est suggérée, demander aux apprenants de répondre à la série de questions dans le cas de la solution suggérée et de vérifier que chaque réponse peut être expliquée."

$r21 = $d.Content
$r21.Find.Execute("Discuss what information they can get from the statement of the riddle") | Out-Null
$r21.Text = "Se demander quelles informations ils peuvent obtenir de l’énoncé de l'énigme"

$r22 = $d.Content
$r22.Find.Execute("Figure out which cases can be excluded ") | Out-Null
$r22.Text = "Déterminer quels sont les cas pouvant être exclus "

$r23 = $d.Content
$r23.Find.Execute("Enact possible solutions") | Out-Null
$r23.Text = "Proposer des solutions possibles"
